$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 55
$ws.Range("I2").Value = 145
$ws.Range("J2").Value = 685
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 201
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 122
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 90
$ws.Range("T2").Value = 119
$ws.Range("U2").Value = 8
$ws.Range("V2").Value = 1049
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1087
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 7
